$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados" timestamp refreshed
$ws.Range("A1").Value = "Datos actualizados a 17 de Octubre de 2020 a las 16:12"

# Full country table refresh (new case counts pulled from source; table re-sorted
# descending by "Casos totales" which also shifts a few countries to new rows).
$data = New-Object 'object[,]' 218,8
$data[0,0] = "Estados Unidos"; $data[0,1] = 8296249; $data[0,2] = 7971; $data[0,3] = 5397672; $data[0,4] = 2674847; $data[0,5] = 0; $data[0,6] = 86; $data[0,7] = 223730
$data[1,0] = "India"; $data[1,1] = 7442823; $data[1,2] = 12188; $data[1,3] = 6533867; $data[1,4] = 795791; $data[1,5] = 0; $data[1,6] = 133; $data[1,7] = 113165
$data[2,0] = "Brasil"; $data[2,1] = 5201570; $data[2,2] = 0; $data[2,3] = 4619560; $data[2,4] = 428781; $data[2,5] = 0; $data[2,6] = 0; $data[2,7] = 153229
$data[3,0] = "Rusia"; $data[3,1] = 1384235; $data[3,2] = 14922; $data[3,3] = 1065199; $data[3,4] = 295034; $data[3,5] = 0; $data[3,6] = 279; $data[3,7] = 24002
$data[4,0] = "España"; $data[4,1] = 982723; $data[4,2] = 0; $data[4,3] = 0; $data[4,4] = 0; $data[4,5] = 0; $data[4,6] = 0; $data[4,7] = 33775
$data[5,0] = "Argentina"; $data[5,1] = 965609; $data[5,2] = 0; $data[5,3] = 778501; $data[5,4] = 161385; $data[5,5] = 0; $data[5,6] = 0; $data[5,7] = 25723
$data[6,0] = "Colombia"; $data[6,1] = 945354; $data[6,2] = 0; $data[6,3] = 837001; $data[6,4] = 79737; $data[6,5] = 0; $data[6,6] = 0; $data[6,7] = 28616
$data[7,0] = "Peru"; $data[7,1] = 862417; $data[7,2] = 0; $data[7,3] = 769077; $data[7,4] = 59692; $data[7,5] = 0; $data[7,6] = 0; $data[7,7] = 33648
$data[8,0] = "Mexico"; $data[8,1] = 841661; $data[8,2] = 6751; $data[8,3] = 612216; $data[8,4] = 143741; $data[8,5] = 0; $data[8,6] = 419; $data[8,7] = 85704
$data[9,0] = "Francia"; $data[9,1] = 834770; $data[9,2] = 0; $data[9,3] = 104696; $data[9,4] = 696771; $data[9,5] = 0; $data[9,6] = 0; $data[9,7] = 33303
$data[10,0] = "Sudafrica"; $data[10,1] = 700203; $data[10,2] = 0; $data[10,3] = 629260; $data[10,4] = 52573; $data[10,5] = 0; $data[10,6] = 0; $data[10,7] = 18370
$data[11,0] = "Reino Unido"; $data[11,1] = 689257; $data[11,2] = 0; $data[11,3] = 0; $data[11,4] = 0; $data[11,5] = 0; $data[11,6] = 0; $data[11,7] = 43429
$data[12,0] = "Iran"; $data[12,1] = 526490; $data[12,2] = 4103; $data[12,3] = 423921; $data[12,4] = 72446; $data[12,5] = 0; $data[12,6] = 253; $data[12,7] = 30123
$data[13,0] = "Chile"; $data[13,1] = 488190; $data[13,2] = 0; $data[13,3] = 461097; $data[13,4] = 13564; $data[13,5] = 0; $data[13,6] = 0; $data[13,7] = 13529
$data[14,0] = "Irak"; $data[14,1] = 423524; $data[14,2] = 3221; $data[14,3] = 357291; $data[14,4] = 56035; $data[14,5] = 0; $data[14,6] = 56; $data[14,7] = 10198
$data[15,0] = "Italia"; $data[15,1] = 391611; $data[15,2] = 0; $data[15,3] = 247872; $data[15,4] = 107312; $data[15,5] = 0; $data[15,6] = 0; $data[15,7] = 36427
$data[16,0] = "Banglades"; $data[16,1] = 387295; $data[16,2] = 1209; $data[16,3] = 302298; $data[16,4] = 79351; $data[16,5] = 0; $data[16,6] = 23; $data[16,7] = 5646
$data[17,0] = "Indonesia"; $data[17,1] = 357762; $data[17,2] = 4301; $data[17,3] = 281592; $data[17,4] = 63739; $data[17,5] = 0; $data[17,6] = 84; $data[17,7] = 12431
$data[18,0] = "Alemania"; $data[18,1] = 357538; $data[18,2] = 746; $data[18,3] = 290000; $data[18,4] = 57700; $data[18,5] = 0; $data[18,6] = 2; $data[18,7] = 9838
$data[19,0] = "Filipinas"; $data[19,1] = 354338; $data[19,2] = 2673; $data[19,3] = 295312; $data[19,4] = 52423; $data[19,5] = 0; $data[19,6] = 73; $data[19,7] = 6603
$data[20,0] = "Turquia"; $data[20,1] = 343955; $data[20,2] = 0; $data[20,3] = 301098; $data[20,4] = 33704; $data[20,5] = 0; $data[20,6] = 0; $data[20,7] = 9153
$data[21,0] = "Arabia Saudita"; $data[21,1] = 341854; $data[21,2] = 359; $data[21,3] = 328165; $data[21,4] = 8524; $data[21,5] = 0; $data[21,6] = 21; $data[21,7] = 5165
$data[22,0] = "Pakistan"; $data[22,1] = 322452; $data[22,2] = 575; $data[22,3] = 306640; $data[22,4] = 9174; $data[22,5] = 0; $data[22,6] = 17; $data[22,7] = 6638
$data[23,0] = "Israel"; $data[23,1] = 301896; $data[23,2] = 0; $data[23,3] = 262503; $data[23,4] = 37252; $data[23,5] = 0; $data[23,6] = 0; $data[23,7] = 2141
$data[24,0] = "Ucrania"; $data[24,1] = 293641; $data[24,2] = 6410; $data[24,3] = 124113; $data[24,4] = 164011; $data[24,5] = 0; $data[24,6] = 109; $data[24,7] = 5517
$data[25,0] = "Paises Bajos"; $data[25,1] = 220052; $data[25,2] = 8114; $data[25,3] = 0; $data[25,4] = 0; $data[25,5] = 0; $data[25,6] = 29; $data[25,7] = 6737
$data[26,0] = "Belgica"; $data[26,1] = 202151; $data[26,2] = 10192; $data[26,3] = 20867; $data[26,4] = 170925; $data[26,5] = 0; $data[26,6] = 32; $data[26,7] = 10359
$data[27,0] = "Canada"; $data[27,1] = 194106; $data[27,2] = 0; $data[27,3] = 163644; $data[27,4] = 20740; $data[27,5] = 0; $data[27,6] = 0; $data[27,7] = 9722
$data[28,0] = "Rumania"; $data[28,1] = 176468; $data[28,2] = 3952; $data[28,3] = 129556; $data[28,4] = 41100; $data[28,5] = 0; $data[28,6] = 63; $data[28,7] = 5812
$data[29,0] = "Polonia"; $data[29,1] = 167230; $data[29,2] = 9622; $data[29,3] = 90162; $data[29,4] = 73544; $data[29,5] = 0; $data[29,6] = 84; $data[29,7] = 3524
$data[30,0] = "Marruecos"; $data[30,1] = 167148; $data[30,2] = 0; $data[30,3] = 138989; $data[30,4] = 25341; $data[30,5] = 0; $data[30,6] = 0; $data[30,7] = 2818
$data[31,0] = "Chequia"; $data[31,1] = 160112; $data[31,2] = 0; $data[31,3] = 66093; $data[31,4] = 92736; $data[31,5] = 0; $data[31,6] = 0; $data[31,7] = 1283
$data[32,0] = "Ecuador"; $data[32,1] = 151659; $data[32,2] = 0; $data[32,3] = 128134; $data[32,4] = 11168; $data[32,5] = 0; $data[32,6] = 0; $data[32,7] = 12357
$data[33,0] = "Bolivia"; $data[33,1] = 139562; $data[33,2] = 243; $data[33,3] = 103896; $data[33,4] = 27227; $data[33,5] = 0; $data[33,6] = 32; $data[33,7] = 8439
$data[34,0] = "Nepal"; $data[34,1] = 129304; $data[34,2] = 3167; $data[34,3] = 89840; $data[34,4] = 38737; $data[34,5] = 0; $data[34,6] = 12; $data[34,7] = 727
$data[35,0] = "Catar"; $data[35,1] = 129227; $data[35,2] = 235; $data[35,3] = 126218; $data[35,4] = 2786; $data[35,5] = 0; $data[35,6] = 1; $data[35,7] = 223
$data[36,0] = "Panama"; $data[36,1] = 123498; $data[36,2] = 0; $data[36,3] = 99286; $data[36,4] = 21666; $data[36,5] = 0; $data[36,6] = 0; $data[36,7] = 2546
$data[37,0] = "Republica Dominicana"; $data[37,1] = 120450; $data[37,2] = 0; $data[37,3] = 96883; $data[37,4] = 21375; $data[37,5] = 0; $data[37,6] = 0; $data[37,7] = 2192
$data[38,0] = "Kuwait"; $data[38,1] = 115483; $data[38,2] = 739; $data[38,3] = 107108; $data[38,4] = 7681; $data[38,5] = 0; $data[38,6] = 4; $data[38,7] = 694
$data[39,0] = "Emiratos Arabes Unidos"; $data[39,1] = 114387; $data[39,2] = 1538; $data[39,3] = 106354; $data[39,4] = 7574; $data[39,5] = 0; $data[39,6] = 4; $data[39,7] = 459
$data[40,0] = "Kazajistan"; $data[40,1] = 109302; $data[40,2] = 100; $data[40,3] = 104921; $data[40,4] = 2613; $data[40,5] = 0; $data[40,6] = 0; $data[40,7] = 1768
$data[41,0] = "Oman"; $data[41,1] = 108296; $data[41,2] = 0; $data[41,3] = 94229; $data[41,4] = 12996; $data[41,5] = 0; $data[41,6] = 0; $data[41,7] = 1071
$data[42,0] = "Egipto"; $data[42,1] = 105159; $data[42,2] = 0; $data[42,3] = 98089; $data[42,4] = 971; $data[42,5] = 0; $data[42,6] = 0; $data[42,7] = 6099
$data[43,0] = "Suecia"; $data[43,1] = 103200; $data[43,2] = 0; $data[43,3] = 0; $data[43,4] = 0; $data[43,5] = 0; $data[43,6] = 0; $data[43,7] = 5918
$data[44,0] = "Guatemala"; $data[44,1] = 100431; $data[44,2] = 0; $data[44,3] = 89494; $data[44,4] = 7459; $data[44,5] = 0; $data[44,6] = 0; $data[44,7] = 3478
$data[45,0] = "Portugal"; $data[45,1] = 98055; $data[45,2] = 2153; $data[45,3] = 57919; $data[45,4] = 37974; $data[45,5] = 0; $data[45,6] = 13; $data[45,7] = 2162
$data[46,0] = "Costa Rica"; $data[46,1] = 94348; $data[46,2] = 0; $data[46,3] = 58269; $data[46,4] = 34911; $data[46,5] = 0; $data[46,6] = 0; $data[46,7] = 1168
$data[47,0] = "Japon"; $data[47,1] = 91431; $data[47,2] = 0; $data[47,3] = 84451; $data[47,4] = 5330; $data[47,5] = 0; $data[47,6] = 0; $data[47,7] = 1650
$data[48,0] = "Etiopia"; $data[48,1] = 87834; $data[48,2] = 0; $data[48,3] = 41628; $data[48,4] = 44869; $data[48,5] = 0; $data[48,6] = 0; $data[48,7] = 1337
$data[49,0] = "Honduras"; $data[49,1] = 86691; $data[49,2] = 602; $data[49,3] = 34546; $data[49,4] = 49589; $data[49,5] = 0; $data[49,6] = 4; $data[49,7] = 2556
$data[50,0] = "Bielorrusia"; $data[50,1] = 86392; $data[50,2] = 0; $data[50,3] = 78990; $data[50,4] = 6481; $data[50,5] = 0; $data[50,6] = 0; $data[50,7] = 921
$data[51,0] = "Venezuela"; $data[51,1] = 85758; $data[51,2] = 0; $data[51,3] = 78294; $data[51,4] = 6739; $data[51,5] = 0; $data[51,6] = 0; $data[51,7] = 725
$data[52,0] = "China"; $data[52,1] = 85659; $data[52,2] = 13; $data[52,3] = 80766; $data[52,4] = 259; $data[52,5] = 0; $data[52,6] = 0; $data[52,7] = 4634
$data[53,0] = "Barein"; $data[53,1] = 77325; $data[53,2] = 0; $data[53,3] = 73421; $data[53,4] = 3611; $data[53,5] = 0; $data[53,6] = 1; $data[53,7] = 293
$data[54,0] = "Suiza"; $data[54,1] = 74422; $data[54,2] = 0; $data[54,3] = 50500; $data[54,4] = 21800; $data[54,5] = 0; $data[54,6] = 0; $data[54,7] = 2122
$data[55,0] = "Moldavia"; $data[55,1] = 65860; $data[55,2] = 0; $data[55,3] = 46543; $data[55,4] = 17768; $data[55,5] = 0; $data[55,6] = 0; $data[55,7] = 1549
$data[56,0] = "Austria"; $data[56,1] = 63134; $data[56,2] = 1747; $data[56,3] = 48771; $data[56,4] = 13474; $data[56,5] = 0; $data[56,6] = 7; $data[56,7] = 889
$data[57,0] = "Armenia"; $data[57,1] = 63000; $data[57,2] = 1540; $data[57,3] = 47925; $data[57,4] = 14008; $data[57,5] = 0; $data[57,6] = 11; $data[57,7] = 1067
$data[58,0] = "Uzbekistan"; $data[58,1] = 62809; $data[58,2] = 221; $data[58,3] = 59756; $data[58,4] = 2531; $data[58,5] = 0; $data[58,6] = 2; $data[58,7] = 522
$data[59,0] = "Nigeria"; $data[59,1] = 61194; $data[59,2] = 0; $data[59,3] = 52304; $data[59,4] = 7771; $data[59,5] = 0; $data[59,6] = 0; $data[59,7] = 1119
$data[60,0] = "Libano"; $data[60,1] = 60113; $data[60,2] = 0; $data[60,3] = 26468; $data[60,4] = 33136; $data[60,5] = 0; $data[60,6] = 0; $data[60,7] = 509
$data[61,0] = "Singapur"; $data[61,1] = 57904; $data[61,2] = 3; $data[61,3] = 57784; $data[61,4] = 92; $data[61,5] = 0; $data[61,6] = 0; $data[61,7] = 28
$data[62,0] = "Argelia"; $data[62,1] = 53998; $data[62,2] = 0; $data[62,3] = 37856; $data[62,4] = 14301; $data[62,5] = 0; $data[62,6] = 0; $data[62,7] = 1841
$data[63,0] = "Paraguay"; $data[63,1] = 53482; $data[63,2] = 0; $data[63,3] = 34927; $data[63,4] = 17390; $data[63,5] = 0; $data[63,6] = 0; $data[63,7] = 1165
$data[64,0] = "Kirguistan"; $data[64,1] = 51490; $data[64,2] = 470; $data[64,3] = 45509; $data[64,4] = 4873; $data[64,5] = 0; $data[64,6] = 5; $data[64,7] = 1108
$data[65,0] = "Libia"; $data[65,1] = 47845; $data[65,2] = 0; $data[65,3] = 26062; $data[65,4] = 21084; $data[65,5] = 0; $data[65,6] = 0; $data[65,7] = 699
$data[66,0] = "Irlanda"; $data[66,1] = 47427; $data[66,2] = 0; $data[66,3] = 23364; $data[66,4] = 22222; $data[66,5] = 0; $data[66,6] = 0; $data[66,7] = 1841
$data[67,0] = "Ghana"; $data[67,1] = 47173; $data[67,2] = 0; $data[67,3] = 46527; $data[67,4] = 336; $data[67,5] = 0; $data[67,6] = 0; $data[67,7] = 310
$data[68,0] = "Estado de Palestina"; $data[68,1] = 46746; $data[68,2] = 312; $data[68,3] = 40162; $data[68,4] = 6182; $data[68,5] = 0; $data[68,6] = 0; $data[68,7] = 402
$data[69,0] = "Hungria"; $data[69,1] = 44816; $data[69,2] = 1791; $data[69,3] = 13580; $data[69,4] = 30127; $data[69,5] = 0; $data[69,6] = 24; $data[69,7] = 1109
$data[70,0] = "Azerbaiyan"; $data[70,1] = 44317; $data[70,2] = 528; $data[70,3] = 39903; $data[70,4] = 3791; $data[70,5] = 0; $data[70,6] = 2; $data[70,7] = 623
$data[71,0] = "Kenia"; $data[71,1] = 44196; $data[71,2] = 616; $data[71,3] = 31752; $data[71,4] = 11619; $data[71,5] = 0; $data[71,6] = 12; $data[71,7] = 825
$data[72,0] = "Afganistan"; $data[72,1] = 40141; $data[72,2] = 68; $data[72,3] = 33561; $data[72,4] = 5092; $data[72,5] = 0; $data[72,6] = 3; $data[72,7] = 1488
$data[73,0] = "Serbia"; $data[73,1] = 35946; $data[73,2] = 227; $data[73,3] = 31536; $data[73,4] = 3636; $data[73,5] = 0; $data[73,6] = 2; $data[73,7] = 774
$data[74,0] = "Dinamarca"; $data[74,1] = 34941; $data[74,2] = 500; $data[74,3] = 28917; $data[74,4] = 5345; $data[74,5] = 0; $data[74,6] = 2; $data[74,7] = 679
$data[75,0] = "Birmania"; $data[75,1] = 34875; $data[75,2] = 1387; $data[75,3] = 16370; $data[75,4] = 17667; $data[75,5] = 0; $data[75,6] = 39; $data[75,7] = 838
$data[76,0] = "Tunez"; $data[76,1] = 34790; $data[76,2] = 0; $data[76,3] = 5032; $data[76,4] = 29246; $data[76,5] = 0; $data[76,6] = 0; $data[76,7] = 512
$data[77,0] = "Jordania"; $data[77,1] = 34548; $data[77,2] = 0; $data[77,3] = 6692; $data[77,4] = 27546; $data[77,5] = 0; $data[77,6] = 0; $data[77,7] = 310
$data[78,0] = "Bosnia y Herzegovina"; $data[78,1] = 33561; $data[78,2] = 716; $data[78,3] = 24773; $data[78,4] = 7807; $data[78,5] = 0; $data[78,6] = 1; $data[78,7] = 981
$data[79,0] = "El Salvador"; $data[79,1] = 31456; $data[79,2] = 191; $data[79,3] = 26769; $data[79,4] = 3770; $data[79,5] = 0; $data[79,6] = 5; $data[79,7] = 917
$data[80,0] = "Bulgaria"; $data[80,1] = 28505; $data[80,2] = 0; $data[80,3] = 16875; $data[80,4] = 10672; $data[80,5] = 0; $data[80,6] = 0; $data[80,7] = 958
$data[81,0] = "Eslovaquia"; $data[81,1] = 28268; $data[81,2] = 1968; $data[81,3] = 7297; $data[81,4] = 20889; $data[81,5] = 0; $data[81,6] = 11; $data[81,7] = 82
$data[82,0] = "Australia"; $data[82,1] = 27383; $data[82,2] = 12; $data[82,3] = 25098; $data[82,4] = 1381; $data[82,5] = 0; $data[82,6] = 0; $data[82,7] = 904
$data[83,0] = "Corea del Sur"; $data[83,1] = 25108; $data[83,2] = 73; $data[83,3] = 23258; $data[83,4] = 1407; $data[83,5] = 0; $data[83,6] = 2; $data[83,7] = 443
$data[84,0] = "Croacia"; $data[84,1] = 24761; $data[84,2] = 1096; $data[84,3] = 19562; $data[84,4] = 4844; $data[84,5] = 0; $data[84,6] = 10; $data[84,7] = 355
$data[85,0] = "Grecia"; $data[85,1] = 24450; $data[85,2] = 0; $data[85,3] = 9989; $data[85,4] = 13971; $data[85,5] = 0; $data[85,6] = 0; $data[85,7] = 490
$data[86,0] = "Republica de Macedonia"; $data[86,1] = 22607; $data[86,2] = 0; $data[86,3] = 16949; $data[86,4] = 4837; $data[86,5] = 0; $data[86,6] = 0; $data[86,7] = 821
$data[87,0] = "Camerun"; $data[87,1] = 21441; $data[87,2] = 0; $data[87,3] = 20117; $data[87,4] = 901; $data[87,5] = 0; $data[87,6] = 0; $data[87,7] = 423
$data[88,0] = "Costa de Marfil"; $data[88,1] = 20275; $data[88,2] = 0; $data[88,3] = 19953; $data[88,4] = 201; $data[88,5] = 0; $data[88,6] = 0; $data[88,7] = 121
$data[89,0] = "Malasia"; $data[89,1] = 19627; $data[89,2] = 869; $data[89,3] = 12561; $data[89,4] = 6886; $data[89,5] = 0; $data[89,6] = 4; $data[89,7] = 180
$data[90,0] = "Madagascar"; $data[90,1] = 16810; $data[90,2] = 56; $data[90,3] = 16215; $data[90,4] = 357; $data[90,5] = 0; $data[90,6] = 1; $data[90,7] = 238
$data[91,0] = "Albania"; $data[91,1] = 16501; $data[91,2] = 0; $data[91,3] = 9957; $data[91,4] = 6101; $data[91,5] = 0; $data[91,6] = 0; $data[91,7] = 443
$data[92,0] = "Noruega"; $data[92,1] = 16349; $data[92,2] = 77; $data[92,3] = 11863; $data[92,4] = 4208; $data[92,5] = 0; $data[92,6] = 0; $data[92,7] = 278
$data[93,0] = "Georgia"; $data[93,1] = 16285; $data[93,2] = 958; $data[93,3] = 7827; $data[93,4] = 8330; $data[93,5] = 0; $data[93,6] = 4; $data[93,7] = 128
$data[94,0] = "Zambia"; $data[94,1] = 15659; $data[94,2] = 0; $data[94,3] = 14899; $data[94,4] = 414; $data[94,5] = 0; $data[94,6] = 0; $data[94,7] = 346
$data[95,0] = "Senegal"; $data[95,1] = 15392; $data[95,2] = 24; $data[95,3] = 13756; $data[95,4] = 1319; $data[95,5] = 0; $data[95,6] = 0; $data[95,7] = 317
$data[96,0] = "Montenegro"; $data[96,1] = 15281; $data[96,2] = 0; $data[96,3] = 10569; $data[96,4] = 4484; $data[96,5] = 0; $data[96,6] = 0; $data[96,7] = 228
$data[97,0] = "Sudan"; $data[97,1] = 13691; $data[97,2] = 0; $data[97,3] = 6764; $data[97,4] = 6091; $data[97,5] = 0; $data[97,6] = 0; $data[97,7] = 836
$data[98,0] = "Finlandia"; $data[98,1] = 13293; $data[98,2] = 160; $data[98,3] = 9100; $data[98,4] = 3842; $data[98,5] = 0; $data[98,6] = 0; $data[98,7] = 351
$data[99,0] = "Eslovenia"; $data[99,1] = 12416; $data[99,2] = 898; $data[99,3] = 6148; $data[99,4] = 6084; $data[99,5] = 0; $data[99,6] = 4; $data[99,7] = 184
$data[100,0] = "Namibia"; $data[100,1] = 12215; $data[100,2] = 0; $data[100,3] = 10360; $data[100,4] = 1724; $data[100,5] = 0; $data[100,6] = 0; $data[100,7] = 131
$data[101,0] = "Guinea"; $data[101,1] = 11362; $data[101,2] = 0; $data[101,3] = 10420; $data[101,4] = 872; $data[101,5] = 0; $data[101,6] = 0; $data[101,7] = 70
$data[102,0] = "Maldivas"; $data[102,1] = 11154; $data[102,2] = 0; $data[102,3] = 9995; $data[102,4] = 1123; $data[102,5] = 0; $data[102,6] = 1; $data[102,7] = 36
$data[103,0] = "Consejo Danes para los Refugiados"; $data[103,1] = 11000; $data[103,2] = 1; $data[103,3] = 10342; $data[103,4] = 356; $data[103,5] = 0; $data[103,6] = 1; $data[103,7] = 302
$data[104,0] = "Mozambique"; $data[104,1] = 10612; $data[104,2] = 0; $data[104,3] = 8262; $data[104,4] = 2277; $data[104,5] = 0; $data[104,6] = 0; $data[104,7] = 73
$data[105,0] = "Luxemburgo"; $data[105,1] = 10471; $data[105,2] = 0; $data[105,3] = 8468; $data[105,4] = 1870; $data[105,5] = 0; $data[105,6] = 0; $data[105,7] = 133
$data[106,0] = "Uganda"; $data[106,1] = 10455; $data[106,2] = 121; $data[106,3] = 6901; $data[106,4] = 3458; $data[106,5] = 0; $data[106,6] = 0; $data[106,7] = 96
$data[107,0] = "Tayikistan"; $data[107,1] = 10455; $data[107,2] = 41; $data[107,3] = 9457; $data[107,4] = 918; $data[107,5] = 0; $data[107,6] = 0; $data[107,7] = 80
$data[108,0] = "Guayana Francesa"; $data[108,1] = 10239; $data[108,2] = 0; $data[108,3] = 9955; $data[108,4] = 215; $data[108,5] = 0; $data[108,6] = 0; $data[108,7] = 69
$data[109,0] = "Haiti"; $data[109,1] = 8925; $data[109,2] = 0; $data[109,3] = 7182; $data[109,4] = 1512; $data[109,5] = 0; $data[109,6] = 0; $data[109,7] = 231
$data[110,0] = "Gabon"; $data[110,1] = 8881; $data[110,2] = 0; $data[110,3] = 8430; $data[110,4] = 397; $data[110,5] = 0; $data[110,6] = 0; $data[110,7] = 54
$data[111,0] = "Jamaica"; $data[111,1] = 8132; $data[111,2] = 0; $data[111,3] = 3653; $data[111,4] = 4317; $data[111,5] = 0; $data[111,6] = 0; $data[111,7] = 162
$data[112,0] = "Zimbabue"; $data[112,1] = 8099; $data[112,2] = 0; $data[112,3] = 7673; $data[112,4] = 195; $data[112,5] = 0; $data[112,6] = 0; $data[112,7] = 231
$data[113,0] = "Mauritania"; $data[113,1] = 7603; $data[113,2] = 0; $data[113,3] = 7339; $data[113,4] = 101; $data[113,5] = 0; $data[113,6] = 0; $data[113,7] = 163
$data[114,0] = "Cabo Verde"; $data[114,1] = 7526; $data[114,2] = 0; $data[114,3] = 6425; $data[114,4] = 1019; $data[114,5] = 0; $data[114,6] = 0; $data[114,7] = 82
$data[115,0] = "Lituania"; $data[115,1] = 7269; $data[115,2] = 228; $data[115,3] = 3097; $data[115,4] = 4059; $data[115,5] = 0; $data[115,6] = 1; $data[115,7] = 113
$data[116,0] = "Angola"; $data[116,1] = 7222; $data[116,2] = 0; $data[116,3] = 3012; $data[116,4] = 3976; $data[116,5] = 0; $data[116,6] = 0; $data[116,7] = 234
$data[117,0] = "Guadalupe"; $data[117,1] = 7122; $data[117,2] = 0; $data[117,3] = 2199; $data[117,4] = 4827; $data[117,5] = 0; $data[117,6] = 0; $data[117,7] = 96
$data[118,0] = "Cuba"; $data[118,1] = 6118; $data[118,2] = 0; $data[118,3] = 5702; $data[118,4] = 292; $data[118,5] = 0; $data[118,6] = 0; $data[118,7] = 124
$data[119,0] = "Malaui"; $data[119,1] = 5842; $data[119,2] = 0; $data[119,3] = 4735; $data[119,4] = 926; $data[119,5] = 0; $data[119,6] = 0; $data[119,7] = 181
$data[120,0] = "Suazilandia"; $data[120,1] = 5746; $data[120,2] = 0; $data[120,3] = 5392; $data[120,4] = 239; $data[120,5] = 0; $data[120,6] = 0; $data[120,7] = 115
$data[121,0] = "Bahamas"; $data[121,1] = 5517; $data[121,2] = 0; $data[121,3] = 3201; $data[121,4] = 2202; $data[121,5] = 0; $data[121,6] = 0; $data[121,7] = 114
$data[122,0] = "Sri Lanka"; $data[122,1] = 5475; $data[122,2] = 121; $data[122,3] = 3395; $data[122,4] = 2067; $data[122,5] = 0; $data[122,6] = 0; $data[122,7] = 13
$data[123,0] = "Republica de Yibuti"; $data[123,1] = 5449; $data[123,2] = 0; $data[123,3] = 5372; $data[123,4] = 16; $data[123,5] = 0; $data[123,6] = 0; $data[123,7] = 61
$data[124,0] = "Nicaragua"; $data[124,1] = 5353; $data[124,2] = 0; $data[124,3] = 4225; $data[124,4] = 974; $data[124,5] = 0; $data[124,6] = 0; $data[124,7] = 154
$data[125,0] = "Botsuana"; $data[125,1] = 5242; $data[125,2] = 0; $data[125,3] = 905; $data[125,4] = 4317; $data[125,5] = 0; $data[125,6] = 0; $data[125,7] = 20
$data[126,0] = "Trinidad yTobago"; $data[126,1] = 5241; $data[126,2] = 0; $data[126,3] = 3545; $data[126,4] = 1601; $data[126,5] = 0; $data[126,6] = 0; $data[126,7] = 95
$data[127,0] = "Hong Kong"; $data[127,1] = 5238; $data[127,2] = 17; $data[127,3] = 4963; $data[127,4] = 170; $data[127,5] = 0; $data[127,6] = 0; $data[127,7] = 105
$data[128,0] = "Congo"; $data[128,1] = 5156; $data[128,2] = 0; $data[128,3] = 3887; $data[128,4] = 1177; $data[128,5] = 0; $data[128,6] = 0; $data[128,7] = 92
$data[129,0] = "Surinam"; $data[129,1] = 5113; $data[129,2] = 0; $data[129,3] = 4921; $data[129,4] = 83; $data[129,5] = 0; $data[129,6] = 0; $data[129,7] = 109
$data[130,0] = "Guinea Ecuatorial"; $data[130,1] = 5068; $data[130,2] = 0; $data[130,3] = 4954; $data[130,4] = 31; $data[130,5] = 0; $data[130,6] = 0; $data[130,7] = 83
$data[131,0] = "Siria"; $data[131,1] = 4987; $data[131,2] = 0; $data[131,3] = 1456; $data[131,4] = 3290; $data[131,5] = 0; $data[131,6] = 0; $data[131,7] = 241
$data[132,0] = "Ruanda"; $data[132,1] = 4965; $data[132,2] = 0; $data[132,3] = 4664; $data[132,4] = 267; $data[132,5] = 0; $data[132,6] = 0; $data[132,7] = 34
$data[133,0] = "Republica de Africa Central"; $data[133,1] = 4855; $data[133,2] = 0; $data[133,3] = 1924; $data[133,4] = 2869; $data[133,5] = 0; $data[133,6] = 0; $data[133,7] = 62
$data[134,0] = "Reunion"; $data[134,1] = 4776; $data[134,2] = 0; $data[134,3] = 4445; $data[134,4] = 314; $data[134,5] = 0; $data[134,6] = 0; $data[134,7] = 17
$data[135,0] = "Malta"; $data[135,1] = 4486; $data[135,2] = 204; $data[135,3] = 3184; $data[135,4] = 1257; $data[135,5] = 0; $data[135,6] = 0; $data[135,7] = 45
$data[136,0] = "Aruba"; $data[136,1] = 4289; $data[136,2] = 0; $data[136,3] = 3947; $data[136,4] = 310; $data[136,5] = 0; $data[136,6] = 0; $data[136,7] = 32
$data[137,0] = "Estonia"; $data[137,1] = 4052; $data[137,2] = 35; $data[137,3] = 3198; $data[137,4] = 786; $data[137,5] = 0; $data[137,6] = 0; $data[137,7] = 68
$data[138,0] = "Mayotte"; $data[138,1] = 4030; $data[138,2] = 0; $data[138,3] = 2964; $data[138,4] = 1023; $data[138,5] = 0; $data[138,6] = 0; $data[138,7] = 43
$data[139,0] = "Islandia"; $data[139,1] = 3998; $data[139,2] = 69; $data[139,3] = 2745; $data[139,4] = 1242; $data[139,5] = 0; $data[139,6] = 0; $data[139,7] = 11
$data[140,0] = "Somalia"; $data[140,1] = 3864; $data[140,2] = 0; $data[140,3] = 3089; $data[140,4] = 676; $data[140,5] = 0; $data[140,6] = 0; $data[140,7] = 99
$data[141,0] = "Polinesia Francesa"; $data[141,1] = 3797; $data[141,2] = 0; $data[141,3] = 2844; $data[141,4] = 939; $data[141,5] = 0; $data[141,6] = 0; $data[141,7] = 14
$data[142,0] = "Tailandia"; $data[142,1] = 3679; $data[142,2] = 10; $data[142,3] = 3478; $data[142,4] = 142; $data[142,5] = 0; $data[142,6] = 0; $data[142,7] = 59
$data[143,0] = "Guyana"; $data[143,1] = 3672; $data[143,2] = 0; $data[143,3] = 2590; $data[143,4] = 975; $data[143,5] = 0; $data[143,6] = 0; $data[143,7] = 107
$data[144,0] = "Gambia"; $data[144,1] = 3649; $data[144,2] = 0; $data[144,3] = 2649; $data[144,4] = 882; $data[144,5] = 0; $data[144,6] = 0; $data[144,7] = 118
$data[145,0] = "Letonia"; $data[145,1] = 3392; $data[145,2] = 188; $data[145,3] = 1329; $data[145,4] = 2020; $data[145,5] = 0; $data[145,6] = 1; $data[145,7] = 43
$data[146,0] = "Mali"; $data[146,1] = 3378; $data[146,2] = 0; $data[146,3] = 2563; $data[146,4] = 683; $data[146,5] = 0; $data[146,6] = 0; $data[146,7] = 132
$data[147,0] = "Principado de Andorra"; $data[147,1] = 3377; $data[147,2] = 0; $data[147,3] = 2057; $data[147,4] = 1261; $data[147,5] = 0; $data[147,6] = 0; $data[147,7] = 59
$data[148,0] = "Sudan del Sur"; $data[148,1] = 2817; $data[148,2] = 0; $data[148,3] = 1290; $data[148,4] = 1472; $data[148,5] = 0; $data[148,6] = 0; $data[148,7] = 55
$data[149,0] = "Belice"; $data[149,1] = 2728; $data[149,2] = 46; $data[149,3] = 1626; $data[149,4] = 1059; $data[149,5] = 0; $data[149,6] = 2; $data[149,7] = 43
$data[150,0] = "Benin"; $data[150,1] = 2496; $data[150,2] = 0; $data[150,3] = 2330; $data[150,4] = 125; $data[150,5] = 0; $data[150,6] = 0; $data[150,7] = 41
$data[151,0] = "Uruguay"; $data[151,1] = 2450; $data[151,2] = 0; $data[151,3] = 2042; $data[151,4] = 357; $data[151,5] = 0; $data[151,6] = 0; $data[151,7] = 51
$data[152,0] = "Guinea-Bisau"; $data[152,1] = 2389; $data[152,2] = 0; $data[152,3] = 1782; $data[152,4] = 566; $data[152,5] = 0; $data[152,6] = 0; $data[152,7] = 41
$data[153,0] = "Republica de Chipre"; $data[153,1] = 2379; $data[153,2] = 0; $data[153,3] = 1444; $data[153,4] = 910; $data[153,5] = 0; $data[153,6] = 0; $data[153,7] = 25
$data[154,0] = "Burkina Faso"; $data[154,1] = 2343; $data[154,2] = 0; $data[154,3] = 1718; $data[154,4] = 560; $data[154,5] = 0; $data[154,6] = 0; $data[154,7] = 65
$data[155,0] = "Sierra Leona"; $data[155,1] = 2325; $data[155,2] = 0; $data[155,3] = 1750; $data[155,4] = 502; $data[155,5] = 0; $data[155,6] = 0; $data[155,7] = 73
$data[156,0] = "Martinica"; $data[156,1] = 2257; $data[156,2] = 0; $data[156,3] = 98; $data[156,4] = 2135; $data[156,5] = 0; $data[156,6] = 0; $data[156,7] = 24
$data[157,0] = "Yemen"; $data[157,1] = 2055; $data[157,2] = 0; $data[157,3] = 1335; $data[157,4] = 124; $data[157,5] = 0; $data[157,6] = 0; $data[157,7] = 596
$data[158,0] = "Togo"; $data[158,1] = 2027; $data[158,2] = 0; $data[158,3] = 1500; $data[158,4] = 476; $data[158,5] = 0; $data[158,6] = 0; $data[158,7] = 51
$data[159,0] = "Nueva Zelanda"; $data[159,1] = 1883; $data[159,2] = 3; $data[159,3] = 1818; $data[159,4] = 40; $data[159,5] = 0; $data[159,6] = 0; $data[159,7] = 25
$data[160,0] = "Lesoto"; $data[160,1] = 1833; $data[160,2] = 0; $data[160,3] = 961; $data[160,4] = 830; $data[160,5] = 0; $data[160,6] = 0; $data[160,7] = 42
$data[161,0] = "Liberia"; $data[161,1] = 1377; $data[161,2] = 0; $data[161,3] = 1264; $data[161,4] = 31; $data[161,5] = 0; $data[161,6] = 0; $data[161,7] = 82
$data[162,0] = "Republica del Chad"; $data[162,1] = 1361; $data[162,2] = 0; $data[162,3] = 1138; $data[162,4] = 130; $data[162,5] = 0; $data[162,6] = 0; $data[162,7] = 93
$data[163,0] = "Niger"; $data[163,1] = 1209; $data[163,2] = 0; $data[163,3] = 1126; $data[163,4] = 14; $data[163,5] = 0; $data[163,6] = 0; $data[163,7] = 69
$data[164,0] = "Vietnam"; $data[164,1] = 1126; $data[164,2] = 2; $data[164,3] = 1031; $data[164,4] = 60; $data[164,5] = 0; $data[164,6] = 0; $data[164,7] = 35
$data[165,0] = "Santo Tome y Principe"; $data[165,1] = 932; $data[165,2] = 0; $data[165,3] = 896; $data[165,4] = 21; $data[165,5] = 0; $data[165,6] = 0; $data[165,7] = 15
$data[166,0] = "San Marino"; $data[166,1] = 759; $data[166,2] = 0; $data[166,3] = 685; $data[166,4] = 32; $data[166,5] = 0; $data[166,6] = 0; $data[166,7] = 42
$data[167,0] = "San Martin (Parte Holandesa)"; $data[167,1] = 746; $data[167,2] = 0; $data[167,3] = 659; $data[167,4] = 65; $data[167,5] = 0; $data[167,6] = 0; $data[167,7] = 22
$data[168,0] = "Crucero"; $data[168,1] = 712; $data[168,2] = 0; $data[168,3] = 659; $data[168,4] = 40; $data[168,5] = 0; $data[168,6] = 0; $data[168,7] = 13
$data[169,0] = "Curazao"; $data[169,1] = 698; $data[169,2] = 0; $data[169,3] = 391; $data[169,4] = 306; $data[169,5] = 0; $data[169,6] = 0; $data[169,7] = 1
$data[170,0] = "Islas Turcas y Caicos"; $data[170,1] = 697; $data[170,2] = 0; $data[170,3] = 674; $data[170,4] = 17; $data[170,5] = 0; $data[170,6] = 0; $data[170,7] = 6
$data[171,0] = "Papua Nueva Guinea"; $data[171,1] = 581; $data[171,2] = 3; $data[171,3] = 540; $data[171,4] = 34; $data[171,5] = 0; $data[171,6] = 0; $data[171,7] = 7
$data[172,0] = "Gibraltar"; $data[172,1] = 558; $data[172,2] = 14; $data[172,3] = 449; $data[172,4] = 109; $data[172,5] = 0; $data[172,6] = 0; $data[172,7] = 0
$data[173,0] = "Taiwan"; $data[173,1] = 535; $data[173,2] = 0; $data[173,3] = 491; $data[173,4] = 37; $data[173,5] = 0; $data[173,6] = 0; $data[173,7] = 7
$data[174,0] = "San Martin (Parte Francesa)"; $data[174,1] = 531; $data[174,2] = 0; $data[174,3] = 380; $data[174,4] = 143; $data[174,5] = 0; $data[174,6] = 0; $data[174,7] = 8
$data[175,0] = "Burundi"; $data[175,1] = 531; $data[175,2] = 0; $data[175,3] = 497; $data[175,4] = 33; $data[175,5] = 0; $data[175,6] = 0; $data[175,7] = 1
$data[176,0] = "Tanzania"; $data[176,1] = 509; $data[176,2] = 0; $data[176,3] = 183; $data[176,4] = 305; $data[176,5] = 0; $data[176,6] = 0; $data[176,7] = 21
$data[177,0] = "Comoras"; $data[177,1] = 502; $data[177,2] = 0; $data[177,3] = 485; $data[177,4] = 10; $data[177,5] = 0; $data[177,6] = 0; $data[177,7] = 7
$data[178,0] = "Islas Feroe"; $data[178,1] = 483; $data[178,2] = 1; $data[178,3] = 472; $data[178,4] = 11; $data[178,5] = 0; $data[178,6] = 0; $data[178,7] = 0
$data[179,0] = "Eritrea"; $data[179,1] = 422; $data[179,2] = 0; $data[179,3] = 376; $data[179,4] = 46; $data[179,5] = 0; $data[179,6] = 0; $data[179,7] = 0
$data[180,0] = "Mauricio"; $data[180,1] = 417; $data[180,2] = 0; $data[180,3] = 364; $data[180,4] = 43; $data[180,5] = 0; $data[180,6] = 0; $data[180,7] = 10
$data[181,0] = "Isla de Man"; $data[181,1] = 348; $data[181,2] = 0; $data[181,3] = 319; $data[181,4] = 5; $data[181,5] = 0; $data[181,6] = 0; $data[181,7] = 24
$data[182,0] = "Mongolia"; $data[182,1] = 320; $data[182,2] = 0; $data[182,3] = 311; $data[182,4] = 9; $data[182,5] = 0; $data[182,6] = 0; $data[182,7] = 0
$data[183,0] = "Butan"; $data[183,1] = 316; $data[183,2] = 0; $data[183,3] = 298; $data[183,4] = 18; $data[183,5] = 0; $data[183,6] = 0; $data[183,7] = 0
$data[184,0] = "Camboya"; $data[184,1] = 283; $data[184,2] = 0; $data[184,3] = 280; $data[184,4] = 3; $data[184,5] = 0; $data[184,6] = 0; $data[184,7] = 0
$data[185,0] = "Monaco"; $data[185,1] = 255; $data[185,2] = 0; $data[185,3] = 217; $data[185,4] = 36; $data[185,5] = 0; $data[185,6] = 0; $data[185,7] = 2
$data[186,0] = "Islas Caimanes"; $data[186,1] = 233; $data[186,2] = 0; $data[186,3] = 212; $data[186,4] = 20; $data[186,5] = 0; $data[186,6] = 0; $data[186,7] = 1
$data[187,0] = "Barbados"; $data[187,1] = 219; $data[187,2] = 0; $data[187,3] = 195; $data[187,4] = 17; $data[187,5] = 0; $data[187,6] = 0; $data[187,7] = 7
$data[188,0] = "Liechtenstein"; $data[188,1] = 192; $data[188,2] = 0; $data[188,3] = 132; $data[188,4] = 59; $data[188,5] = 0; $data[188,6] = 0; $data[188,7] = 1
$data[189,0] = "Bermudas"; $data[189,1] = 185; $data[189,2] = 0; $data[189,3] = 172; $data[189,4] = 4; $data[189,5] = 0; $data[189,6] = 0; $data[189,7] = 9
$data[190,0] = "Bonaire, San Eustaquio y Saba"; $data[190,1] = 150; $data[190,2] = 0; $data[190,3] = 111; $data[190,4] = 37; $data[190,5] = 0; $data[190,6] = 0; $data[190,7] = 2
$data[191,0] = "Seychelles"; $data[191,1] = 149; $data[191,2] = 0; $data[191,3] = 148; $data[191,4] = 1; $data[191,5] = 0; $data[191,6] = 0; $data[191,7] = 0
$data[192,0] = "Brunei"; $data[192,1] = 147; $data[192,2] = 0; $data[192,3] = 143; $data[192,4] = 1; $data[192,5] = 0; $data[192,6] = 0; $data[192,7] = 3
$data[193,0] = "Antigua y Barbuda"; $data[193,1] = 112; $data[193,2] = 0; $data[193,3] = 100; $data[193,4] = 9; $data[193,5] = 0; $data[193,6] = 0; $data[193,7] = 3
$data[194,0] = "San Bartolome"; $data[194,1] = 72; $data[194,2] = 0; $data[194,3] = 55; $data[194,4] = 17; $data[194,5] = 0; $data[194,6] = 0; $data[194,7] = 0
$data[195,0] = "Islas Virgenes Britanicas"; $data[195,1] = 71; $data[195,2] = 0; $data[195,3] = 70; $data[195,4] = 0; $data[195,5] = 0; $data[195,6] = 0; $data[195,7] = 1
$data[196,0] = "San Vicente y las Granadinas"; $data[196,1] = 65; $data[196,2] = 0; $data[196,3] = 64; $data[196,4] = 1; $data[196,5] = 0; $data[196,6] = 0; $data[196,7] = 0
$data[197,0] = "Macao"; $data[197,1] = 46; $data[197,2] = 0; $data[197,3] = 46; $data[197,4] = 0; $data[197,5] = 0; $data[197,6] = 0; $data[197,7] = 0
$data[198,0] = "Puerto Rico"; $data[198,1] = 39; $data[198,2] = 0; $data[198,3] = 1; $data[198,4] = 36; $data[198,5] = 0; $data[198,6] = 0; $data[198,7] = 2
$data[199,0] = "Dominica"; $data[199,1] = 33; $data[199,2] = 0; $data[199,3] = 29; $data[199,4] = 4; $data[199,5] = 0; $data[199,6] = 0; $data[199,7] = 0
$data[200,0] = "Guam"; $data[200,1] = 32; $data[200,2] = 0; $data[200,3] = 0; $data[200,4] = 31; $data[200,5] = 0; $data[200,6] = 0; $data[200,7] = 1
$data[201,0] = "Santa Lucia"; $data[201,1] = 32; $data[201,2] = 0; $data[201,3] = 27; $data[201,4] = 5; $data[201,5] = 0; $data[201,6] = 0; $data[201,7] = 0
$data[202,0] = "Fiyi"; $data[202,1] = 32; $data[202,2] = 0; $data[202,3] = 30; $data[202,4] = 0; $data[202,5] = 0; $data[202,6] = 0; $data[202,7] = 2
$data[203,0] = "Timor Oriental"; $data[203,1] = 29; $data[203,2] = 0; $data[203,3] = 28; $data[203,4] = 1; $data[203,5] = 0; $data[203,6] = 0; $data[203,7] = 0
$data[204,0] = "Santa Sede"; $data[204,1] = 27; $data[204,2] = 1; $data[204,3] = 15; $data[204,4] = 12; $data[204,5] = 0; $data[204,6] = 0; $data[204,7] = 0
$data[205,0] = "Nueva Caledonia"; $data[205,1] = 27; $data[205,2] = 0; $data[205,3] = 27; $data[205,4] = 0; $data[205,5] = 0; $data[205,6] = 0; $data[205,7] = 0
$data[206,0] = "Granada"; $data[206,1] = 25; $data[206,2] = 0; $data[206,3] = 24; $data[206,4] = 1; $data[206,5] = 0; $data[206,6] = 0; $data[206,7] = 0
$data[207,0] = "Laos"; $data[207,1] = 23; $data[207,2] = 0; $data[207,3] = 22; $data[207,4] = 1; $data[207,5] = 0; $data[207,6] = 0; $data[207,7] = 0
$data[208,0] = "San Cristobal y Nieves"; $data[208,1] = 19; $data[208,2] = 0; $data[208,3] = 19; $data[208,4] = 0; $data[208,5] = 0; $data[208,6] = 0; $data[208,7] = 0
$data[209,0] = "Islas Virgenes de los Estados Unidos"; $data[209,1] = 17; $data[209,2] = 0; $data[209,3] = 0; $data[209,4] = 17; $data[209,5] = 0; $data[209,6] = 0; $data[209,7] = 0
$data[210,0] = "San Pedro y Miquelon"; $data[210,1] = 16; $data[210,2] = 0; $data[210,3] = 12; $data[210,4] = 4; $data[210,5] = 0; $data[210,6] = 0; $data[210,7] = 0
$data[211,0] = "Groenlandia"; $data[211,1] = 16; $data[211,2] = 0; $data[211,3] = 14; $data[211,4] = 2; $data[211,5] = 0; $data[211,6] = 0; $data[211,7] = 0
$data[212,0] = "Montserrat"; $data[212,1] = 13; $data[212,2] = 0; $data[212,3] = 12; $data[212,4] = 0; $data[212,5] = 0; $data[212,6] = 0; $data[212,7] = 1
$data[213,0] = "Islas Malvinas"; $data[213,1] = 13; $data[213,2] = 0; $data[213,3] = 13; $data[213,4] = 0; $data[213,5] = 0; $data[213,6] = 0; $data[213,7] = 0
$data[214,0] = "Sahara Occidental"; $data[214,1] = 10; $data[214,2] = 0; $data[214,3] = 8; $data[214,4] = 1; $data[214,5] = 0; $data[214,6] = 0; $data[214,7] = 1
$data[215,0] = "Islas Salomon"; $data[215,1] = 3; $data[215,2] = 0; $data[215,3] = 0; $data[215,4] = 3; $data[215,5] = 0; $data[215,6] = 0; $data[215,7] = 0
$data[216,0] = "Anguila"; $data[216,1] = 3; $data[216,2] = 0; $data[216,3] = 3; $data[216,4] = 0; $data[216,5] = 0; $data[216,6] = 0; $data[216,7] = 0
$data[217,0] = "Wallis y Futuna"; $data[217,1] = 1; $data[217,2] = 0; $data[217,3] = 0; $data[217,4] = 1; $data[217,5] = 0; $data[217,6] = 0; $data[217,7] = 0

$ws.Range("A4:H221").Value = $data
